$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B (Book Title) - closest achievable value to the target 54.14
$ws.Columns.Item(2).ColumnWidth = 53.3

# Publisher Code column (M) - renumber BP-000025 -> BP-00001 for every data row
$ws.Range("M2:M17").Value = "BP-00001"

# Authors column (R) - renumber / restructure author codes per row
$ws.Range("R2").Value = "A-000001"
$ws.Range("R3").Value = "A-000002"
$ws.Range("R4").Value = "A-000001,A-000002,A-000003"
$ws.Range("R5").Value = "A-000001"
$ws.Range("R6").Value = "A-000001"
$ws.Range("R7").Value = "A-000004,A-000005,A-000006"
$ws.Range("R8").Value = "A-000001"
$ws.Range("R9").Value = "A-000001"
$ws.Range("R10").Value = "A-000001"
$ws.Range("R11").Value = "A-000001"
$ws.Range("R12").Value = "A-000001"
$ws.Range("R13").Value = "A-000001"
$ws.Range("R14").Value = "A-000001"
$ws.Range("R15").Value = "A-000001"
$ws.Range("R16").Value = "A-000001"
$ws.Range("R17").Value = "A-000001"
